# New installer generation for version 1.2
# Refresh the ground-station contact predictions (rows 1-5: header + the
# first four passes) with the newly regenerated ephemeris data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (text content unchanged, just re-asserted) ---
$ws.Range("A1").Value = "Source"
$ws.Range("B1").Value = "Target"
$ws.Range("C1").Value = "IntervalNumber"
$ws.Range("D1").Value = "StartTime"
$ws.Range("E1").Value = "EndTime"
$ws.Range("F1").Value = "Duration"
$ws.Range("G1").Value = "StartOrbit"
$ws.Range("H1").Value = "EndOrbit"

# --- Row 2 (interval 1) ---
$ws.Range("A2").Value = "experiment 3"
$ws.Range("B2").Value = "ESA Kiruna Ground Station Rx"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 47844.479861111111
$ws.Range("E2").Value = 47844.512499999997
$ws.Range("F2").Value = 2820
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1

# --- Row 3 (interval 2) ---
$ws.Range("A3").Value = "experiment 3"
$ws.Range("B3").Value = "ESA Kiruna Ground Station Rx"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 47845.478125000001
$ws.Range("E3").Value = 47845.511111111111
$ws.Range("F3").Value = 2850
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 3

# --- Row 4 (interval 3) ---
$ws.Range("A4").Value = "experiment 3"
$ws.Range("B4").Value = "ESA Kiruna Ground Station Rx"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 47846.476388888892
$ws.Range("E4").Value = 47846.509722222225
$ws.Range("F4").Value = 2880
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 5

# --- Row 5 (interval 4) ---
$ws.Range("A5").Value = "experiment 3"
$ws.Range("B5").Value = "ESA Kiruna Ground Station Rx"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 47847.474652777775
$ws.Range("E5").Value = 47847.508333333331
$ws.Range("F5").Value = 2910
$ws.Range("G5").Value = 7
$ws.Range("H5").Value = 7
